$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-05-26 Sunday" "2024-05-27 Monday"

Replace-Text "827×3=2481" "345×4=1380"
Replace-Text "581×6=3486" "154×3=462"
Replace-Text "247×4=988" "859×7=6013"
Replace-Text "749×8=5992" "223×8=1784"
Replace-Text "717×5=3585" "760×2=1520"

Replace-Text "768×8=6144" "163×6=978"
Replace-Text "633×4=2532" "226×9=2034"
Replace-Text "371×9=3339" "535×8=4280"
Replace-Text "918×8=7344" "676×7=4732"
Replace-Text "591×9=5319" "174×2=348"

Replace-Text "466×6=2796" "453×8=3624"
Replace-Text "330×2=660" "172×9=1548"
Replace-Text "907×4=3628" "300×8=2400"
Replace-Text "935×6=5610" "280×6=1680"
Replace-Text "456×2=912" "536×7=3752"

Replace-Text "878×5=4390" "236×3=708"
Replace-Text "581×2=1162" "505×9=4545"
Replace-Text "115×7=805" "728×8=5824"
Replace-Text "557×8=4456" "481×8=3848"
Replace-Text "650×3=1950" "564×4=2256"

Replace-Text "689×8=5512" "735×9=6615"
Replace-Text "306×9=2754" "410×9=3690"
Replace-Text "876×5=4380" "578×4=2312"
Replace-Text "416×3=1248" "196×4=784"
Replace-Text "301×8=2408" "494×8=3952"
